$wb = $excel.ActiveWorkbook

# ----- Sheet: LP1912 -----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 14:19:39"
$ws.Cells.Item(3,1).Value = "Total filas: 274"
$ws.Cells.Item(39,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(40,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(65,1).Value = "05:52:07"
$ws.Cells.Item(65,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(65,4).Value = 100
$ws.Cells.Item(66,1).Value = "07:28:14"
$ws.Cells.Item(66,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(66,4).Value = 4
$ws.Cells.Item(72,1).Value = "06:59:37"
$ws.Cells.Item(72,3).Value = "14_ABASTO"
$ws.Cells.Item(72,4).Value = 48
$ws.Cells.Item(73,1).Value = "07:28:14"
$ws.Cells.Item(73,3).Value = "16_SANTA ANA"
$ws.Cells.Item(73,4).Value = 19
$ws.Cells.Item(75,3).Value = "10_OLMOS"
$ws.Cells.Item(76,3).Value = "215D_EL PATO"
$ws.Cells.Item(88,1).Value = "07:28:14"
$ws.Cells.Item(88,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(88,4).Value = 55
$ws.Cells.Item(89,1).Value = "08:13:38"
$ws.Cells.Item(89,3).Value = "215B_EL PATO"
$ws.Cells.Item(89,4).Value = 10
$ws.Cells.Item(114,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(115,3).Value = "16_SANTA ANA"
$ws.Cells.Item(118,1).Value = "09:33:33"
$ws.Cells.Item(118,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(118,4).Value = 0
$ws.Cells.Item(119,1).Value = "08:56:26"
$ws.Cells.Item(119,3).Value = "10_OLMOS"
$ws.Cells.Item(119,4).Value = 37
$ws.Cells.Item(153,3).Value = "27_EL RETIRO"
$ws.Cells.Item(154,3).Value = "16_SANTA ANA"
$ws.Cells.Item(193,1).Value = "12:02:21"
$ws.Cells.Item(193,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(193,4).Value = 32
$ws.Cells.Item(194,1).Value = "11:49:23"
$ws.Cells.Item(194,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(194,4).Value = 45
$ws.Cells.Item(199,3).Value = "17_179 Y 38"
$ws.Cells.Item(200,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(226,3).Value = "14_ABASTO"
$ws.Cells.Item(227,3).Value = "15_ABASTO"
$ws.Cells.Item(236,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(237,3).Value = "215A_EL PATO"
$ws.Cells.Item(243,1).Value = "13:18:32"
$ws.Cells.Item(243,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(243,4).Value = 46
$ws.Cells.Item(244,1).Value = "13:51:20"
$ws.Cells.Item(244,3).Value = "17_ROMERO"
$ws.Cells.Item(244,4).Value = 13
$ws.Cells.Item(250,1).Value = "14:19:39"
$ws.Cells.Item(250,2).Value = "14:19"
$ws.Cells.Item(250,4).Value = 0
$ws.Cells.Item(251,2).Value = "14:20"
$ws.Cells.Item(251,3).Value = "215C_EL PATO"
$ws.Cells.Item(251,4).Value = 29
$ws.Cells.Item(252,1).Value = "14:19:39"
$ws.Cells.Item(252,2).Value = "14:21"
$ws.Cells.Item(252,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(252,4).Value = 2
$ws.Cells.Item(253,1).Value = "14:19:39"
$ws.Cells.Item(253,2).Value = "14:22"
$ws.Cells.Item(253,4).Value = 3
$ws.Cells.Item(254,1).Value = "14:19:39"
$ws.Cells.Item(254,2).Value = "14:28"
$ws.Cells.Item(254,3).Value = "15_ABASTO"
$ws.Cells.Item(254,4).Value = 9
$ws.Cells.Item(255,2).Value = "14:30"
$ws.Cells.Item(255,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(255,4).Value = 39
$ws.Cells.Item(256,1).Value = "14:19:39"
$ws.Cells.Item(256,2).Value = "14:34"
$ws.Cells.Item(256,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(256,4).Value = 15
$ws.Cells.Item(257,1).Value = "12:54:24"
$ws.Cells.Item(257,2).Value = "14:39"
$ws.Cells.Item(257,3).Value = "14_ABASTO"
$ws.Cells.Item(257,4).Value = 105
$ws.Cells.Item(258,1).Value = "14:19:39"
$ws.Cells.Item(258,2).Value = "14:44"
$ws.Cells.Item(258,3).Value = "14_ABASTO"
$ws.Cells.Item(258,4).Value = 25
$ws.Cells.Item(259,1).Value = "14:19:39"
$ws.Cells.Item(259,2).Value = "14:46"
$ws.Cells.Item(259,3).Value = "16_SANTA ANA"
$ws.Cells.Item(259,4).Value = 27
$ws.Cells.Item(260,1).Value = "14:19:39"
$ws.Cells.Item(260,2).Value = "14:56"
$ws.Cells.Item(260,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(260,4).Value = 37
$ws.Cells.Item(261,1).Value = "14:19:39"
$ws.Cells.Item(261,2).Value = "14:58"
$ws.Cells.Item(261,3).Value = "215B_EL PATO"
$ws.Cells.Item(261,4).Value = 39
$ws.Cells.Item(262,1).Value = "14:19:39"
$ws.Cells.Item(262,2).Value = "15:00"
$ws.Cells.Item(262,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(262,4).Value = 41
$ws.Cells.Item(263,1).Value = "14:19:39"
$ws.Cells.Item(263,2).Value = "15:05"
$ws.Cells.Item(263,3).Value = "10_OLMOS"
$ws.Cells.Item(263,4).Value = 46
$ws.Cells.Item(264,1).Value = "14:19:39"
$ws.Cells.Item(264,2).Value = "15:06"
$ws.Cells.Item(264,3).Value = "16_SANTA ANA"
$ws.Cells.Item(264,4).Value = 47
$ws.Cells.Item(265,1).Value = "14:19:39"
$ws.Cells.Item(265,2).Value = "15:10"
$ws.Cells.Item(265,3).Value = "17_ROMERO"
$ws.Cells.Item(265,4).Value = 51
$ws.Cells.Item(266,1).Value = "14:19:39"
$ws.Cells.Item(266,2).Value = "15:13"
$ws.Cells.Item(266,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(266,4).Value = 54
$ws.Cells.Item(267,1).Value = "14:19:39"
$ws.Cells.Item(267,2).Value = "15:20"
$ws.Cells.Item(267,3).Value = "15_ABASTO"
$ws.Cells.Item(267,4).Value = 61
$ws.Cells.Item(268,1).Value = "14:19:39"
$ws.Cells.Item(268,2).Value = "15:21"
$ws.Cells.Item(268,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(268,4).Value = 62
$ws.Cells.Item(269,1).Value = "14:19:39"
$ws.Cells.Item(269,2).Value = "15:32"
$ws.Cells.Item(269,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(269,4).Value = 73
$ws.Cells.Item(269,5).Value = "LP1912"
$ws.Cells.Item(270,1).Value = "14:19:39"
$ws.Cells.Item(270,2).Value = "15:34"
$ws.Cells.Item(270,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(270,4).Value = 75
$ws.Cells.Item(270,5).Value = "LP1912"
$ws.Cells.Item(271,1).Value = "13:51:20"
$ws.Cells.Item(271,2).Value = "15:35"
$ws.Cells.Item(271,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(271,4).Value = 104
$ws.Cells.Item(271,5).Value = "LP1912"
$ws.Cells.Item(272,1).Value = "14:19:39"
$ws.Cells.Item(272,2).Value = "15:37"
$ws.Cells.Item(272,3).Value = "10_OLMOS"
$ws.Cells.Item(272,4).Value = 78
$ws.Cells.Item(272,5).Value = "LP1912"
$ws.Cells.Item(273,1).Value = "14:19:39"
$ws.Cells.Item(273,2).Value = "15:38"
$ws.Cells.Item(273,3).Value = "215A_EL PATO"
$ws.Cells.Item(273,4).Value = 79
$ws.Cells.Item(273,5).Value = "LP1912"
$ws.Cells.Item(274,1).Value = "13:51:20"
$ws.Cells.Item(274,2).Value = "15:44"
$ws.Cells.Item(274,3).Value = "14_ABASTO"
$ws.Cells.Item(274,4).Value = 113
$ws.Cells.Item(274,5).Value = "LP1912"
$ws.Cells.Item(275,1).Value = "14:19:39"
$ws.Cells.Item(275,2).Value = "15:46"
$ws.Cells.Item(275,3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(275,4).Value = 87
$ws.Cells.Item(275,5).Value = "LP1912"
$ws.Cells.Item(276,1).Value = "14:19:39"
$ws.Cells.Item(276,2).Value = "15:53"
$ws.Cells.Item(276,3).Value = "27_EL RETIRO"
$ws.Cells.Item(276,4).Value = 94
$ws.Cells.Item(276,5).Value = "LP1912"
$ws.Cells.Item(277,1).Value = "14:19:39"
$ws.Cells.Item(277,2).Value = "15:53"
$ws.Cells.Item(277,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(277,4).Value = 94
$ws.Cells.Item(277,5).Value = "LP1912"
$ws.Cells.Item(278,1).Value = "14:19:39"
$ws.Cells.Item(278,2).Value = "15:55"
$ws.Cells.Item(278,3).Value = "17_ROMERO"
$ws.Cells.Item(278,4).Value = 96
$ws.Cells.Item(278,5).Value = "LP1912"
$ws.Cells.Item(279,1).Value = "14:19:39"
$ws.Cells.Item(279,2).Value = "16:15"
$ws.Cells.Item(279,3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(279,4).Value = 116
$ws.Cells.Item(279,5).Value = "LP1912"

# ----- Sheet: LP1912-215 -----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 14:19:39"
$ws.Cells.Item(3,1).Value = "Total filas: 33"
$ws.Cells.Item(35,1).Value = "14:19:39"
$ws.Cells.Item(35,2).Value = "14:19"
$ws.Cells.Item(35,4).Value = 0
$ws.Cells.Item(36,2).Value = "14:20"
$ws.Cells.Item(36,3).Value = "215C_EL PATO"
$ws.Cells.Item(36,4).Value = 29
$ws.Cells.Item(37,1).Value = "14:19:39"
$ws.Cells.Item(37,2).Value = "14:58"
$ws.Cells.Item(37,3).Value = "215B_EL PATO"
$ws.Cells.Item(37,4).Value = 39
$ws.Cells.Item(38,1).Value = "14:19:39"
$ws.Cells.Item(38,2).Value = "15:38"
$ws.Cells.Item(38,3).Value = "215A_EL PATO"
$ws.Cells.Item(38,4).Value = 79
$ws.Cells.Item(38,5).Value = "LP1912"

# ----- Sheet: 6203-6173 -----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 14:19:39"
$ws.Cells.Item(3,1).Value = "Total filas: 39"
$ws.Cells.Item(40,1).Value = "14:19:39"
$ws.Cells.Item(40,2).Value = "14:52"
$ws.Cells.Item(40,4).Value = 33
$ws.Cells.Item(41,2).Value = "14:53"
$ws.Cells.Item(41,3).Value = "215D_LA PLATA"
$ws.Cells.Item(41,4).Value = 62
$ws.Cells.Item(41,5).Value = "L6203"
$ws.Cells.Item(42,1).Value = "13:51:20"
$ws.Cells.Item(42,2).Value = "15:34"
$ws.Cells.Item(42,3).Value = "215A_LA PLATA"
$ws.Cells.Item(42,4).Value = 103
$ws.Cells.Item(42,5).Value = "L6173"
$ws.Cells.Item(43,1).Value = "14:19:39"
$ws.Cells.Item(43,2).Value = "15:35"
$ws.Cells.Item(43,3).Value = "215A_LA PLATA"
$ws.Cells.Item(43,4).Value = 76
$ws.Cells.Item(43,5).Value = "L6173"
$ws.Cells.Item(44,1).Value = "14:19:39"
$ws.Cells.Item(44,2).Value = "16:13"
$ws.Cells.Item(44,3).Value = "215C_LA PLATA"
$ws.Cells.Item(44,4).Value = 114
$ws.Cells.Item(44,5).Value = "L6203"
